# "Generate Report for Handoff"
#
# The localization status report is refreshed: file b.md has now been
# handed off for localization (zh-cn and de-de). Its previous handback
# file is stale relative to the freshly generated handoff, so an error
# detail is recorded for it, and its "Content Duplicate" flag is cleared.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fa8987de776f04a1ce75eadb0e196dd698ffdd28/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/264011c8059e4f3f5247eb05be3b6b5f10a24985/e2e/b.md."

# Helper: write a literal text value even when it looks like a boolean
# keyword ("True"/"False") so Excel doesn't silently coerce it to a
# Boolean cell - a leading apostrophe forces text entry, then resetting
# the style back to Normal clears the quote-prefix flag again.
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# --- Overview sheet: row for b.md ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-25 20:35:51"

# --- zh-cn sheet: row 3 is b.md ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
Set-TextValue $wsZhCn.Range("F3") "False"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-25 20:35:46"
$wsZhCn.Range("P3").Value = $errorDetail

# --- de-de sheet: row 3 is b.md ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
Set-TextValue $wsDeDe.Range("F3") "False"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-25 20:35:51"
$wsDeDe.Range("P3").Value = $errorDetail

# Widen the Error Detail column on both language sheets so the long
# message is readable.
$wsZhCn.Columns.Item(16).ColumnWidth = 40
$wsDeDe.Columns.Item(16).ColumnWidth = 40
